# The "OrganizationType" staging template table gained two columns
# ("Code" and "Name") in the regenerated source table, inserted so the
# header row reads: OrganizationType_ID, BusinessKey, Code, Description, Name
# (i.e. "Code" pushed in ahead of the pre-existing "Description" column,
# and "Name" appended after it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing header cell C2 used to read "Description" - it now reads "Code",
# and the former "Description" column moves out to D2, with a brand new
# "Name" column added at E2.
$ws.Range("C2").Value = "Code"
$ws.Range("D2").Value = "Description"
$ws.Range("E2").Value = "Name"
